$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force Text format first so the literal string (incl. trailing zeros) is kept.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.56"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.44"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.134"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.89"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.410"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000203"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.36"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.93"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.28"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.73"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "420.99"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.607"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.98"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000118"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.27"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.00"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.48"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.156"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.78"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.40"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.67"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.43"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.64"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.20"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0843"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.15"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.875"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.94"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.83"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.47"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.24"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.14"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.39"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.09"

# Remaining text cells (coin names, links, and price/volume strings that
# are not ambiguous with numbers) can be assigned directly.
$ws.Range("D2").Value = "66.328.39"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "3.569.62"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("D7").Value = "3.566.91"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "4.175.91"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("D16").Value = "3.581.72"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "66.351.34"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  -4.28%  "
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").Value = "3.717.24"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "3.571.90"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("E33").Value = "  +3.29%  "
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("E35").Value = "  -3.57%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("E38").Value = "  -3.12%  "
$ws.Range("E39").Value = "  -4.80%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  -4.96%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("E49").Value = "  -4.93%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("E51").Value = "  -0.77%  "
